$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A - Id
$ws.Range("A2").Value = 102077473

# Column B - Taxonsorteringsordning
$ws.Range("B2").Value = 96367

# Column E - TaxonId
$ws.Range("E2").Value = 219874

# Column F - Artnamn
$ws.Range("F2").Value = "Nattviol"

# Column G - Vetenskapligt namn
$ws.Range("G2").Value = "Platanthera bifolia"

# Column H - Auktor
$ws.Range("H2").Value = "(L.) Rich."

# Column I - Antal
$ws.Range("I2").Value = "3"

# Columns J, K, L, N - newly present but empty
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("N2").Value = ""

# Column P - Lokalnamn
$ws.Range("P2").Value = "Tvetaspåret, Tveta, Srm"

# Column Q - Ost
$ws.Range("Q2").Value = 647720.9098417715

# Column R - Nord
$ws.Range("R2").Value = 6560694.968483768

# Column S - Noggrannhet
$ws.Range("S2").Value = 10

# Column Y - Startdatum
$ws.Range("Y2").Value = "2022-06-28"

# Column AA - Slutdatum
$ws.Range("AA2").Value = "2022-07-05"

# Column AF - newly present but empty
$ws.Range("AF2").Value = ""

# Column AI - Biotop-beskrivning, cell removed/cleared entirely
$ws.Range("AI2").ClearContents()

# Column AW - Rapportör
$ws.Range("AW2").Value = "Åsa Johansson"

# Column AX - Observatörer
$ws.Range("AX2").Value = "Åsa Johansson"
